$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 5) to the summary sheet, mirroring the
# existing HFC rows with a new date label ("12/20/16") and the same
# Frequency / Cumulative Frequency / Percent Target / Cumulative Percent
# Target values as the prior rows (0, 20, 0, 20).

# Column A holds text-formatted dates (e.g. "11/21/16"), not real date
# serials, so force the cell to text before entry and then clear the
# temporary formatting back to the default so it stays a plain string.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "12/20/16"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 20
